# Auto-generated: apply scheduled-runner price/profit refresh to Sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3000
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()

$ws.Range("H80").Value = 2655.6
$ws.Range("I80").Value = 506.4
$ws.Range("J80").Value = 3730.2
$ws.Range("K80").Value = 1519.2
$ws.Range("L80").Value = 11190.6
$ws.Range("M80").Value = -521.1999999999998
$ws.Range("N80").Value = -13186.6

$ws.Range("H83").Value = 2655.6
$ws.Range("I83").Value = 506.4
$ws.Range("J83").Value = 3730.2
$ws.Range("K83").Value = 4557.599999999999
$ws.Range("L83").Value = 33571.8
$ws.Range("M83").Value = 434.4000000000005
$ws.Range("N83").Value = -43555.8

$ws.Range("H88").Value = 5474.625
$ws.Range("I88").Value = 3931.3333
$ws.Range("J88").Value = 6400.6
$ws.Range("K88").Value = 3931.3333
$ws.Range("L88").Value = 6400.6
$ws.Range("M88").Value = -3525.3333
$ws.Range("N88").Value = -7212.6

$ws.Range("H91").Value = 5474.625
$ws.Range("I91").Value = 3931.3333
$ws.Range("J91").Value = 6400.6
$ws.Range("K91").Value = 3931.3333
$ws.Range("L91").Value = 6400.6
$ws.Range("M91").Value = -2527.3333
$ws.Range("N91").Value = -9208.6

$ws.Range("H128").Value = 33375
$ws.Range("J128").Value = 33375
$ws.Range("L128").Value = 33375
$ws.Range("N128").Value = -43335

$ws.Range("H129").Value = 2506.682
$ws.Range("I129").Value = 2174
$ws.Range("J129").Value = 2696.7856
$ws.Range("K129").Value = 6522
$ws.Range("L129").Value = 8090.3568
$ws.Range("M129").Value = -1522
$ws.Range("N129").Value = -18090.3568

$ws.Range("H137").Value = 3503.6155
$ws.Range("I137").Value = 2955.1667
$ws.Range("J137").Value = 3973.7144
$ws.Range("K137").Value = 8865.500100000001
$ws.Range("L137").Value = 11921.1432
$ws.Range("M137").Value = -6315.500100000001
$ws.Range("N137").Value = -17021.1432

$ws.Range("H141").Value = 4802.4375
$ws.Range("I141").Value = 3902.8333
$ws.Range("J141").Value = 7501.25
$ws.Range("K141").Value = 11708.4999
$ws.Range("L141").Value = 22503.75
$ws.Range("M141").Value = -6528.499899999999
$ws.Range("N141").Value = -32863.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 326.63635
$ws.Range("I5").Value = 130
$ws.Range("K5").Value = 130
$ws.Range("M5").Value = -18

$ws.Range("H32").Value = 2361.423
$ws.Range("I32").Value = 1527.5682
$ws.Range("K32").Value = 1527.5682
$ws.Range("M32").Value = -1240.5682

$ws.Range("H102").Value = 1736
$ws.Range("I102").Value = 1829.25
$ws.Range("K102").Value = 1829.25
$ws.Range("M102").Value = -207.25

$ws.Range("H132").Value = 8831
$ws.Range("I132").Value = 3904.5
$ws.Range("J132").Value = 10801.6
$ws.Range("K132").Value = 11713.5
$ws.Range("L132").Value = 32404.8
$ws.Range("M132").Value = -9183.5
$ws.Range("N132").Value = -37464.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 326.63635
$ws.Range("I4").Value = 130
$ws.Range("K4").Value = 130
$ws.Range("M4").Value = -15

$ws.Range("H8").Value = 19999.334
$ws.Range("I8").Value = 25999
$ws.Range("J8").Value = 8000
$ws.Range("K8").Value = 25999
$ws.Range("L8").Value = 8000
$ws.Range("M8").Value = -25859
$ws.Range("N8").Value = -8280

$ws.Range("H99").Value = 2816.6667
$ws.Range("I99").Value = 2466.6667
$ws.Range("J99").Value = 3166.6667
$ws.Range("K99").Value = 2466.6667
$ws.Range("L99").Value = 3166.6667
$ws.Range("M99").Value = -968.6667000000002
$ws.Range("N99").Value = -6162.6667

$ws.Range("H107").Value = 1475.375
$ws.Range("J107").Value = 2916.6667
$ws.Range("L107").Value = 2916.6667
$ws.Range("N107").Value = -6756.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 135.33333
$ws.Range("J7").Value = 229.25
$ws.Range("L7").Value = 229.25
$ws.Range("N7").Value = -455.25

$ws.Range("H31").Value = 36489.727
$ws.Range("J31").Value = 71313.375
$ws.Range("L31").Value = 71313.375
$ws.Range("N31").Value = -71903.375

$ws.Range("H34").Value = 36489.727
$ws.Range("J34").Value = 71313.375
$ws.Range("L34").Value = 71313.375
$ws.Range("N34").Value = -71717.375

$ws.Range("H45").Value = 4575
$ws.Range("I45").Value = 5500
$ws.Range("K45").Value = 5500
$ws.Range("M45").Value = -4907

$ws.Range("H58").Value = 6050.893
$ws.Range("I58").Value = 5182.2383
$ws.Range("K58").Value = 5182.2383
$ws.Range("M58").Value = -4979.2383

$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()

$ws.Range("H132").Value = 4454.6343
$ws.Range("I132").Value = 4432.8857
$ws.Range("K132").Value = 13298.6571
$ws.Range("M132").Value = -10768.6571

$ws.Range("H134").Value = 3115.4666
$ws.Range("I134").Value = 2472.9546
$ws.Range("J134").Value = 4882.375
$ws.Range("K134").Value = 7418.8638
$ws.Range("L134").Value = 14647.125
$ws.Range("M134").Value = -4883.8638
$ws.Range("N134").Value = -19717.125

$ws.Range("H136").Value = 6050.893
$ws.Range("I136").Value = 5182.2383
$ws.Range("K136").Value = 15546.7149
$ws.Range("M136").Value = -12996.7149

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 12965.223
$ws.Range("J32").Value = 12965.223
$ws.Range("L32").Value = 38895.669
$ws.Range("N32").Value = -39461.669

$ws.Range("H59").Value = 999
$ws.Range("I59").Value = 999
$ws.Range("K59").Value = 2997
$ws.Range("M59").Value = -2457

$ws.Range("H81").Value = 148092.14
$ws.Range("I81").Value = 2365
$ws.Range("K81").Value = 7095
$ws.Range("M81").Value = -5972

$ws.Range("H84").Value = 148092.14
$ws.Range("I84").Value = 2365
$ws.Range("K84").Value = 21285
$ws.Range("M84").Value = -15669

$ws.Range("H131").Value = 12154346
$ws.Range("J131").Value = 38890090
$ws.Range("L131").Value = 116670270
$ws.Range("N131").Value = -116680350

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 1311.4
$ws.Range("I2").Value = 5.5
$ws.Range("J2").Value = 1637.875
$ws.Range("K2").Value = 5.5
$ws.Range("L2").Value = 1637.875
$ws.Range("M2").Value = 107.5
$ws.Range("N2").Value = -1863.875

$ws.Range("H102").Value = 2024.5518
$ws.Range("I102").Value = 1446
$ws.Range("J102").Value = 4801.6
$ws.Range("K102").Value = 1446
$ws.Range("L102").Value = 4801.6
$ws.Range("M102").Value = 176
$ws.Range("N102").Value = -8045.6

$ws.Range("H122").Value = 3828.0667
$ws.Range("I122").Value = 2108.077
$ws.Range("J122").Value = 15008
$ws.Range("K122").Value = 6324.231000000001
$ws.Range("L122").Value = 45024
$ws.Range("M122").Value = -3874.231000000001
$ws.Range("N122").Value = -49924

$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

$ws.Range("H132").Value = 210802.6
$ws.Range("I132").Value = 999999
$ws.Range("J132").Value = 13503.5
$ws.Range("K132").Value = 2999997
$ws.Range("L132").Value = 40510.5
$ws.Range("M132").Value = -2997467
$ws.Range("N132").Value = -45570.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5190.933
$ws.Range("I7").Value = 2443.762
$ws.Range("J7").Value = 11601
$ws.Range("K7").Value = 2443.762
$ws.Range("L7").Value = 11601
$ws.Range("M7").Value = -2331.762
$ws.Range("N7").Value = -11825

$ws.Range("H40").Value = 3349.0588
$ws.Range("I40").Value = 2370.5625
$ws.Range("J40").Value = 19005
$ws.Range("K40").Value = 2370.5625
$ws.Range("L40").Value = 19005
$ws.Range("M40").Value = -2234.5625
$ws.Range("N40").Value = -19277

$ws.Range("H126").Value = 5190.933
$ws.Range("I126").Value = 2443.762
$ws.Range("J126").Value = 11601
$ws.Range("K126").Value = 7331.286
$ws.Range("L126").Value = 34803
$ws.Range("M126").Value = -4861.286
$ws.Range("N126").Value = -39743

$ws.Range("H136").Value = 4206.5815
$ws.Range("I136").Value = 2144.2903
$ws.Range("J136").Value = 9534.166999999999
$ws.Range("K136").Value = 6432.8709
$ws.Range("L136").Value = 28602.501
$ws.Range("M136").Value = -3882.8709
$ws.Range("N136").Value = -33702.501

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()

$ws.Range("H122").Value = 7148.8823
$ws.Range("I122").Value = 2388.3157
$ws.Range("J122").Value = 13178.934
$ws.Range("K122").Value = 7164.9471
$ws.Range("L122").Value = 39536.802
$ws.Range("M122").Value = -4714.9471
$ws.Range("N122").Value = -44436.802

$ws.Range("H126").Value = 2097.5
$ws.Range("I126").Value = 1157.2222
$ws.Range("J126").Value = 4918.3335
$ws.Range("K126").Value = 3471.6666
$ws.Range("L126").Value = 14755.0005
$ws.Range("M126").Value = -1001.6666
$ws.Range("N126").Value = -19695.0005

$ws.Range("H132").Value = 4828.4814
$ws.Range("I132").Value = 4581.087
$ws.Range("K132").Value = 13743.261
$ws.Range("M132").Value = -11213.261
